$p = $ppt.ActivePresentation
$s1 = $p.Slides.Item(1)
$tbl = $s1.Shapes.Item(2).Table

# Row 11: "Descanso" -> "Estudio Independiente"
$tbl.Cell(11, 1).Shape.TextFrame.TextRange.Text = "Estudio Independiente"

# Row 13: "Estudio Independiente " -> "Mercado "
$tbl.Cell(13, 1).Shape.TextFrame.TextRange.Text = "Mercado "
